# Update "26_mar.xlsx" sheet:
#  - insert two new header rows at the top (shifting the existing table down by 2 rows)
#  - first new row: bold, centered, thin-bordered labels "Unnamed: 0/1/2"
#  - second new row: plain text labels "cidade" / "Casos confirmados" / "Obitos confirmados"
#  - append two new data rows at the bottom: "outros estados" (8) and "outros paises" (28)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing table (rows 1-35) down by two rows.
$ws.Range("A1:A2").EntireRow.Insert()

# New row 1: header-ish labels, bold + centered + thin border all around.
$ws.Cells.Item(1, 1).Value = "Unnamed: 0"
$ws.Cells.Item(1, 2).Value = "Unnamed: 1"
$ws.Cells.Item(1, 3).Value = "Unnamed: 2"

$hdr = $ws.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# New row 2: plain column names.
$ws.Cells.Item(2, 1).Value = "cidade"
$ws.Cells.Item(2, 2).Value = "Casos confirmados"
$ws.Cells.Item(2, 3).Value = "Óbitos confirmados"

# Two brand-new rows appended at the end of the data.
$ws.Cells.Item(38, 1).Value = "outros estados"
$ws.Cells.Item(38, 2).Value = 8

$ws.Cells.Item(39, 1).Value = "outros paises"
$ws.Cells.Item(39, 2).Value = 28
